$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row (row 1) ---
$ws.Range("A1").Value = "...1"
$ws.Range("B1").Value = "rhat_max"
$ws.Range("C1").Value = "ESS_min"
$ws.Range("D1").Value = "pvalue_d1_gr1"
$ws.Range("E1").Value = "pvalue_d1_gr2"
$ws.Range("F1").Value = "pvalue_d2_gr1"
$ws.Range("G1").Value = "pvalue_d2_gr2"
$ws.Range("H1").Value = "pvalue_d3_gr1"
$ws.Range("I1").Value = "pvalue_d3_gr2"
$ws.Range("J1").Value = "waic_d1"
$ws.Range("K1").Value = "waic_d2"
$ws.Range("L1").Value = "waic_d3"
$ws.Range("M1").Value = "waic_tot"
$ws.Range("N1").Value = "CV_d1"
$ws.Range("O1").Value = "CV_d2"
$ws.Range("P1").Value = "CV_d3"
$ws.Range("Q1").Value = "CV_tot"
$ws.Range("R1").Value = "model"
$ws.Range("S1").Value = "beta_log_dist_to_shore"
$ws.Range("T1").Value = "beta_log_bathymetry"
$ws.Range("U1").Value = "beta_mean_CHL"
$ws.Range("V1").Value = "beta_sd_SAL"
$ws.Range("W1").Value = "beta_mean_SSH"
$ws.Range("X1").Value = "sd_beta_log_dist_to_shore"
$ws.Range("Y1").Value = "sd_beta_log_bathymetry"
$ws.Range("Z1").Value = "sd_beta_mean_CHL"
$ws.Range("AA1").Value = "sd_beta_sd_SAL"
$ws.Range("AB1").Value = "sd_beta_mean_SSH"

# --- Update row 2 ---
$ws.Range("A2").Value = "without spatial"
$ws.Range("B2").Value = 1.06
$ws.Range("C2").Value = 321.2
$ws.Range("D2").Value = 0.67
$ws.Range("E2").Value = 0.5
$ws.Range("F2").Value = 0.21
$ws.Range("G2").Value = 0.51
$ws.Range("H2").Value = 0.65
$ws.Range("I2").Value = 0.5
$ws.Range("J2").Value = 750
$ws.Range("K2").Value = 237
$ws.Range("L2").Value = 455
$ws.Range("M2").Value = 1442
$ws.Range("N2").Value = 749
$ws.Range("O2").Value = 232
$ws.Range("P2").Value = 448
$ws.Range("Q2").Value = 1429
$ws.Range("R2").Value = 1
$ws.Range("S2").Value = -1.36
$ws.Range("T2").Value = -0.67
$ws.Range("U2").Value = 2.46
$ws.Range("V2").Value = -0.77
$ws.Range("W2").Value = 0.37
$ws.Range("X2").Value = 0.53
$ws.Range("Y2").Value = 0.35
$ws.Range("Z2").Value = 0.78
$ws.Range("AA2").Value = 0.5
$ws.Range("AB2").Value = 0.42

# --- Update row 3 ---
$ws.Range("A3").Value = "with spatial exp"
$ws.Range("B3").Value = 1.843
$ws.Range("C3").Value = 19.3
$ws.Range("D3").Value = 0.67
$ws.Range("E3").Value = 0.5
$ws.Range("F3").Value = 0.22
$ws.Range("G3").Value = 0.51
$ws.Range("H3").Value = 0.64
$ws.Range("I3").Value = 0.5
$ws.Range("J3").Value = 748
$ws.Range("K3").Value = 237
$ws.Range("L3").Value = 453
$ws.Range("M3").Value = 1438
$ws.Range("N3").Value = 846
$ws.Range("O3").Value = 267
$ws.Range("P3").Value = 494
$ws.Range("Q3").Value = 1607
$ws.Range("R3").Value = 2
$ws.Range("S3").Value = -1.49
$ws.Range("T3").Value = -0.65
$ws.Range("U3").Value = 2.42
$ws.Range("V3").Value = -0.66
$ws.Range("W3").Value = 0.65
$ws.Range("X3").Value = 0.66
$ws.Range("Y3").Value = 0.48
$ws.Range("Z3").Value = 0.94
$ws.Range("AA3").Value = 0.62
$ws.Range("AB3").Value = 0.54

# --- Update row 4 ---
$ws.Range("A4").Value = "with sp shpere"
$ws.Range("B4").Value = 4.014
$ws.Range("C4").Value = 18.9
$ws.Range("D4").Value = 0.67
$ws.Range("E4").Value = 0.5
$ws.Range("F4").Value = 0.21
$ws.Range("G4").Value = 0.5
$ws.Range("H4").Value = 0.63
$ws.Range("I4").Value = 0.5
$ws.Range("J4").Value = 750
$ws.Range("K4").Value = 238
$ws.Range("L4").Value = 454
$ws.Range("M4").Value = 1442
$ws.Range("N4").Value = 846
$ws.Range("O4").Value = 268
$ws.Range("P4").Value = 494
$ws.Range("Q4").Value = 1608
$ws.Range("R4").Value = 3
$ws.Range("S4").Value = -1.53
$ws.Range("T4").Value = -0.75
$ws.Range("U4").Value = 2.46
$ws.Range("V4").Value = -0.69
$ws.Range("W4").Value = 0.66
$ws.Range("X4").Value = 0.66
$ws.Range("Y4").Value = 0.49
$ws.Range("Z4").Value = 0.98
$ws.Range("AA4").Value = 0.7
$ws.Range("AB4").Value = 0.59

# --- Update row 5 ---
$ws.Range("A5").Value = "with sp gaussian"
$ws.Range("B5").Value = 4.991
$ws.Range("C5").Value = 16.9
$ws.Range("D5").Value = 0.67
$ws.Range("E5").Value = 0.5
$ws.Range("F5").Value = 0.22
$ws.Range("G5").Value = 0.5
$ws.Range("H5").Value = 0.64
$ws.Range("I5").Value = 0.5
$ws.Range("J5").Value = 751
$ws.Range("K5").Value = 238
$ws.Range("L5").Value = 455
$ws.Range("M5").Value = 1444
$ws.Range("N5").Value = 845
$ws.Range("O5").Value = 267
$ws.Range("P5").Value = 494
$ws.Range("Q5").Value = 1606
$ws.Range("R5").Value = 4
$ws.Range("S5").Value = -1.43
$ws.Range("T5").Value = -0.73
$ws.Range("U5").Value = 2.49
$ws.Range("V5").Value = -0.81
$ws.Range("W5").Value = 0.44
$ws.Range("X5").Value = 0.6
$ws.Range("Y5").Value = 0.39
$ws.Range("Z5").Value = 0.88
$ws.Range("AA5").Value = 0.56
$ws.Range("AB5").Value = 0.46
# --- Update conditional formatting expression thresholds ---
function Set-ExpressionFormula {
    param($rangeAddr, $newFormula)
    $fcs = $ws.Range($rangeAddr).FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fc = $fcs.Item($i)
        if ($fc.Type -eq 2) {
            $fc.Formula1 = $newFormula
        }
    }
}

Set-ExpressionFormula "J2:J5" "=J2<751"
Set-ExpressionFormula "K2:K5" "=K2<240"
Set-ExpressionFormula "L2:L5" "=L2<456"
Set-ExpressionFormula "M2:M5" "=M2<1441"
Set-ExpressionFormula "N2:N5" "=N2<752"
Set-ExpressionFormula "O2:O5" "=O2<235"
Set-ExpressionFormula "P2:P5" "=P2<451"
Set-ExpressionFormula "Q2:Q5" "=Q2<1432"

Write-Output "All updates applied successfully"
